$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 42605.888090277775
$ws.Range("B3").Value = -30
$ws.Range("C3").Value = 49
$ws.Range("D3").Value = 49
$ws.Range("E3").Value = 9
$ws.Range("F3").Value = 90
$ws.Range("G3").Value = 28684
$ws.Range("H3").Value = 7978
$ws.Range("I3").Value = 420
$ws.Range("J3").Value = 58
$ws.Range("K3").Value = 57
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 10
$ws.Range("N3").Value = "Named"
